# Time Log.xlsx - add the missing time-tracking entry for row 98
# (Sheet1) and move the active-cell selection down to C99, as part of
# the "card layout" fix (issue #27). All downstream formulas (Delta in
# column E, the grand total in E104, the SUMIF/percentage rollups on
# Sheet2, and the pie-chart series) recalc automatically from these
# raw inputs.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# New log entry: 10/23/2014, 10:29 PM -> 11:54 PM, 5 min interruption,
# activity = Coding. Interruption (D98) is written before the Start/Stop
# times so the Delta formula in E98 (which only turns "live" once both
# B98 and C98 are non-blank) picks up the correct Interruption value the
# moment it first evaluates.
$ws1.Range("A98").Value = 41935
$ws1.Range("D98").Value = 5
$ws1.Range("B98").Value = 0.93680555555555556
$ws1.Range("C98").Value = 0.99583333333333324
$ws1.Range("F98").Value = "Coding"

# Make sure every formula (Delta, the SUM total, the Sheet2 SUMIFs,
# and the percentages) is recalculated with the new row in place.
$excel.Calculate()

# The author's selection ended up one row further down afterwards.
$ws1.Range("C99").Select()
